# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# D = Price, E = Volume(1h) change. Values that look like plain numbers (e.g.
# "11.30", "4.32") are written with a leading apostrophe so Excel keeps them
# as text (matching the sheet's existing inlineStr/text storage) instead of
# silently reinterpreting them as numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.583.04'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '1.814.71'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''228.21'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = '''0.559'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '''34.77'
$ws.Range('E8').Value = '  +7.53%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').Value = '''0.0694'
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('D11').Value = '''0.0952'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '2.077.65'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '''11.30'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').Value = '1.818.16'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('E15').Value = '  +2.43%  '
$ws.Range('D16').Value = '34.622.23'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '''4.32'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '''69.10'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').Value = '''247.33'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '0.0₃0801'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').Value = '''11.54'
$ws.Range('E21').Value = '  +5.37%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '''4.22'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('D24').Value = '''172.03'
$ws.Range('E24').Value = '  +5.96%  '
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('D26').Value = '''7.46'
$ws.Range('E26').Value = '  +3.86%  '
$ws.Range('D27').Value = '''16.75'
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''4.04'
$ws.Range('E30').Value = '  +5.55%  '
$ws.Range('D31').Value = '''0.0531'
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').Value = '''3.85'
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').Value = '''1.25'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('E34').Value = '  +2.33%  '
$ws.Range('D35').Value = '''2.61'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = '1.418.17'
$ws.Range('E36').Value = '  -1.73%  '
$ws.Range('D37').Value = '''0.676'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').Value = '''1.07'
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').Value = '''85.86'
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').Value = '''0.962'
$ws.Range('E42').Value = '  +4.06%  '
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = '''13.91'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '''0.0525'
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('E46').Value = '  +3.02%  '
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').Value = '1.978.07'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('D49').Value = '''105.99'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('E51').Value = '  +0.16%  '
